$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.8644080082785512, 0.6696657573405245, 0.7425248273359284, 3.1606896692062696)"
$ws.Range("C2").Value = "JSU(-1.2915927192390875, 1.309534380286125, 1.8199811698037167, 6.101416012982963)"
$ws.Range("D2").Value = "EXN(2.2158431998352235, 0.7959555729601491, 1.8192212856747798)"
$ws.Range("E2").Value = "JSU(-1.6766275585473513, 1.2662011397194, 1.5031266593958772, 4.445613733018956)"
